$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6; this shifts the existing rows 6..108
# down to 7..109, preserving all their data and formatting (matches the
# rest of the diff, which is just every old row N (N>=6) becoming row N+1).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new data record.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 44921
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100103
$ws.Range("H6").Value = "Frutos de hueso (carozo)"
$ws.Range("I6").Value = 100103001
$ws.Range("J6").Value = "Cereza"
$ws.Range("K6").Value = "Brooks"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 800
$ws.Range("N6").Value = 5000
$ws.Range("O6").Value = 5500
$ws.Range("P6").Value = 5250
$ws.Range("Q6").Value = '$/bandeja 10 kilos'
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 525
$ws.Range("T6").Value = 10
